# Applies the row-swap / odds corrections described in the diff
# (league DB refresh commit, 02-05-2024 20:28)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 181
$ws.Range("B181").Value = 6810169
$ws.Range("E181").Value = 'Westerlo'
$ws.Range("F181").Value = 'Cercle Brugge'
$ws.Range("G181").Value = 4
$ws.Range("H181").Value = 2
$ws.Range("I181").Value = 'H'
$ws.Range("J181").Value = 3.4
$ws.Range("K181").Value = 3.8
$ws.Range("L181").Value = 1.95
$ws.Range("M181").Value = 3.3
$ws.Range("N181").Value = 3.6
$ws.Range("O181").Value = 2.05
$ws.Range("P181").Value = 0.25
$ws.Range("Q181").Value = 2
$ws.Range("R181").Value = 1.85
$ws.Range("S181").Value = 2.75
$ws.Range("T181").Value = 1.975
$ws.Range("U181").Value = 1.875
$ws.Range("V181").Value = 2.3
$ws.Range("W181").Value = -1
$ws.Range("Y181").Value = 1
$ws.Range("Z181").Value = -1
$ws.Range("AA181").Value = 0.9750000000000001

# Row 182
$ws.Range("B182").Value = 6810167
$ws.Range("E182").Value = 'Club Brugge'
$ws.Range("F182").Value = 'KV Kortrijk'
$ws.Range("G182").Value = 3
$ws.Range("H182").Value = 3
$ws.Range("I182").Value = 'D'
$ws.Range("J182").Value = 1.125
$ws.Range("K182").Value = 8.5
$ws.Range("L182").Value = 17
$ws.Range("M182").Value = 1.125
$ws.Range("N182").Value = 8.5
$ws.Range("O182").Value = 17
$ws.Range("P182").Value = -2.25
$ws.Range("Q182").Value = 1.85
$ws.Range("R182").Value = 2
$ws.Range("S182").Value = 3.5
$ws.Range("T182").Value = 1.9
$ws.Range("U182").Value = 1.95
$ws.Range("V182").Value = -1
$ws.Range("W182").Value = 7.5
$ws.Range("Y182").Value = -1
$ws.Range("Z182").Value = 1
$ws.Range("AA182").Value = 0.8999999999999999

# Row 183
$ws.Range("B183").Value = 6810168
$ws.Range("E183").Value = 'OH Leuven'
$ws.Range("F183").Value = 'Genk'
$ws.Range("G183").Value = 2
$ws.Range("H183").Value = 1
$ws.Range("J183").Value = 4.5
$ws.Range("K183").Value = 4.2
$ws.Range("L183").Value = 1.666
$ws.Range("M183").Value = 4.333
$ws.Range("N183").Value = 4
$ws.Range("O183").Value = 1.7
$ws.Range("P183").Value = 0.75
$ws.Range("Q183").Value = 1.95
$ws.Range("R183").Value = 1.9
$ws.Range("S183").Value = 3
$ws.Range("T183").Value = 1.975
$ws.Range("U183").Value = 1.875
$ws.Range("V183").Value = 3.333
$ws.Range("Y183").Value = 0.95
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0
$ws.Range("AB183").Value = 0

# Row 184
$ws.Range("B184").Value = 6810165
$ws.Range("E184").Value = 'Charleroi'
$ws.Range("F184").Value = 'Eupen'
$ws.Range("G184").Value = 1
$ws.Range("H184").Value = 0
$ws.Range("J184").Value = 1.6
$ws.Range("K184").Value = 4
$ws.Range("L184").Value = 5
$ws.Range("M184").Value = 1.8
$ws.Range("N184").Value = 3.75
$ws.Range("O184").Value = 4.2
$ws.Range("P184").Value = -0.75
$ws.Range("Q184").Value = 2.05
$ws.Range("R184").Value = 1.8
$ws.Range("S184").Value = 2.75
$ws.Range("T184").Value = 1.95
$ws.Range("U184").Value = 1.9
$ws.Range("V184").Value = 0.8
$ws.Range("Y184").Value = 0.5249999999999999
$ws.Range("Z184").Value = -0.5
$ws.Range("AA184").Value = -1
$ws.Range("AB184").Value = 0.8999999999999999

# Row 187
$ws.Range("B187").Value = 6810163
$ws.Range("E187").Value = 'SintTruidense'
$ws.Range("F187").Value = 'Gent'
$ws.Range("G187").Value = 4
$ws.Range("H187").Value = 1
$ws.Range("I187").Value = 'H'
$ws.Range("J187").Value = 3.6
$ws.Range("K187").Value = 3.6
$ws.Range("L187").Value = 1.95
$ws.Range("M187").Value = 3.25
$ws.Range("N187").Value = 3.4
$ws.Range("O187").Value = 2.15
$ws.Range("Q187").Value = 1.95
$ws.Range("R187").Value = 1.9
$ws.Range("T187").Value = 1.975
$ws.Range("U187").Value = 1.875
$ws.Range("V187").Value = 2.25
$ws.Range("W187").Value = -1
$ws.Range("Y187").Value = 0.95
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 0.9750000000000001

# Row 188
$ws.Range("B188").Value = 6810166
$ws.Range("E188").Value = 'KV Mechelen'
$ws.Range("F188").Value = 'Anderlecht'
$ws.Range("G188").Value = 2
$ws.Range("H188").Value = 2
$ws.Range("I188").Value = 'D'
$ws.Range("J188").Value = 3.5
$ws.Range("K188").Value = 3.5
$ws.Range("L188").Value = 2
$ws.Range("M188").Value = 3
$ws.Range("N188").Value = 3.5
$ws.Range("O188").Value = 2.2
$ws.Range("Q188").Value = 1.925
$ws.Range("R188").Value = 1.925
$ws.Range("T188").Value = 1.875
$ws.Range("U188").Value = 1.975
$ws.Range("V188").Value = -1
$ws.Range("W188").Value = 2.5
$ws.Range("Y188").Value = 0.4625
$ws.Range("Z188").Value = -0.5
$ws.Range("AA188").Value = 0.875

# Row 190
$ws.Range("B190").Value = 6810171
$ws.Range("E190").Value = 'KV Kortrijk'
$ws.Range("F190").Value = 'Charleroi'
$ws.Range("G190").Value = 1
$ws.Range("H190").Value = 0
$ws.Range("I190").Value = 'H'
$ws.Range("J190").Value = 3.2
$ws.Range("K190").Value = 3.5
$ws.Range("L190").Value = 2.1
$ws.Range("M190").Value = 3.4
$ws.Range("N190").Value = 3.4
$ws.Range("O190").Value = 2.05
$ws.Range("P190").Value = 0.25
$ws.Range("Q190").Value = 2
$ws.Range("R190").Value = 1.85
$ws.Range("T190").Value = 1.925
$ws.Range("U190").Value = 1.925
$ws.Range("V190").Value = 2.4
$ws.Range("X190").Value = -1
$ws.Range("Y190").Value = 1
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = 0.925

# Row 191
$ws.Range("B191").Value = 6810174
$ws.Range("E191").Value = 'Westerlo'
$ws.Range("F191").Value = 'OH Leuven'
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 3
$ws.Range("I191").Value = 'A'
$ws.Range("J191").Value = 1.909
$ws.Range("K191").Value = 3.75
$ws.Range("L191").Value = 3.5
$ws.Range("M191").Value = 1.909
$ws.Range("N191").Value = 3.5
$ws.Range("O191").Value = 3.8
$ws.Range("P191").Value = -0.5
$ws.Range("Q191").Value = 1.925
$ws.Range("R191").Value = 1.925
$ws.Range("T191").Value = 1.85
$ws.Range("U191").Value = 2
$ws.Range("V191").Value = -1
$ws.Range("X191").Value = 2.8
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = 0.925
$ws.Range("AA191").Value = 0.8500000000000001
$ws.Range("AB191").Value = -1

# Row 289
$ws.Range("M289").Value = 2.3
$ws.Range("O289").Value = 2.75
$ws.Range("P289").Value = -0.25
$ws.Range("Q289").Value = 2.05
$ws.Range("R289").Value = 1.8

# Row 290
$ws.Range("M290").Value = 2.5
$ws.Range("N290").Value = 3.4
$ws.Range("O290").Value = 2.7
$ws.Range("Q290").Value = 1.825
$ws.Range("R290").Value = 2.025

# Row 291
$ws.Range("P291").Value = 0.5
$ws.Range("Q291").Value = 1.8
$ws.Range("R291").Value = 2.05
$ws.Range("T291").Value = 1.85
$ws.Range("U291").Value = 2

# Row 293
$ws.Range("Q293").Value = 1.9
$ws.Range("R293").Value = 1.95
$ws.Range("S293").Value = 2.5
$ws.Range("T293").Value = 1.8
$ws.Range("U293").Value = 2.05

# Row 295
$ws.Range("M295").Value = 1.7
$ws.Range("N295").Value = 3.8
$ws.Range("Q295").Value = 1.925
$ws.Range("R295").Value = 1.925
$ws.Range("T295").Value = 2
$ws.Range("U295").Value = 1.85

# Row 296
$ws.Range("T296").Value = 1.9
$ws.Range("U296").Value = 1.95
